$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("1031:1031").Delete()
$ws.Rows("1028:1028").Delete()
